$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new detail columns (W, X) to the header row (row 7) and the
# corresponding data row (row 8) of the "Add Panels" sheet.
$ws.Range("W7").Value = "AlarmLoadingDetail"
$ws.Range("X7").Value = "StandbyLoadingDetail"
$ws.Range("W8").Value = "Battery Alarm (A)"
$ws.Range("X8").Value = "Battery Standby (A)"

# Match the formatting of the existing header / data cells in the same rows.
$ws.Range("A7").Copy()
$ws.Range("W7:X7").PasteSpecial(-4122)

$ws.Range("A8").Copy()
$ws.Range("W8:X8").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Reflect the new selection left on the sheet after the edit.
$ws.Range("W7:X8").Select()
